$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: B10/C10/E10 were stored as text ("56348.0", "-2074.0", "0");
# convert them to genuine numeric cells (56348, -2074, 0). D10/F10/G10 stay untouched.
$ws.Cells.Item(10, 2).Value = 56348
$ws.Cells.Item(10, 3).Value = -2074
$ws.Cells.Item(10, 5).Value = 0

# --- Row 11: brand-new row of data, all text-typed (matching the sheet's existing
# convention of storing these report values as inline/shared strings).
# Values that "look" numeric/date/percent get auto-coerced by the COM value setter
# (to a date serial, a number, a percentage...) together with a derived number
# format style; NumberFormat "@" beforehand keeps them literal text, and
# ClearFormats() afterwards drops the incidental style so the cell ends up
# unstyled text, matching the target.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(11, 1) "2022-01-11"
Set-TextValue $ws.Cells.Item(11, 2) "56308.0"
Set-TextValue $ws.Cells.Item(11, 3) "-1656.0"
Set-TextValue $ws.Cells.Item(11, 4) "-2.94%"
Set-TextValue $ws.Cells.Item(11, 5) "-3264.0"

# F11/G11 are present but empty (empty inline/shared string), just like F10/G10.
# A bare "" assignment clears the cell instead of leaving an empty text cell behind,
# so use a leading apostrophe (forces an empty-text literal) and then drop the
# resulting quote-prefix style, same trick as above.
Set-TextValue $ws.Cells.Item(11, 6) "'"
Set-TextValue $ws.Cells.Item(11, 7) "'"
